# Weekly update: insert the newest week's record at row 86, pushing the
# existing historical rows (old 86-154) down by one row (to 87-155).
#
# The data for the other columns of the new record reuse the price
# bracket / origin that were previously on row 86 (only the date and the
# traded volume are new), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 86; this shifts rows 86..154
# down to 87..155 and also bumps the sheet dimension to A1:R155.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the latest week's data.
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C86").Value = "Ñuble"
$ws.Range("D86").Value = 44447
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = 100112008
$ws.Range("G86").Value = "Coliflor"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 160
$ws.Range("K86").Value = 700
$ws.Range("L86").Value = 750
$ws.Range("M86").Value = 725
$ws.Range("N86").Value = "$/unidad"
$ws.Range("O86").Value = "Región del Maule"
$ws.Range("P86").Value = 725
$ws.Range("Q86").Value = 1
$ws.Range("R86").Value = "Hortaliza"
